# "support blank row and colum start"
# Inserts a new "Sheet2" worksheet between Sheet1 and Sheet3, fills it with data,
# and shifts the data on Sheet3 down one row / right one column (leaving a blank
# row 1 and blank column A), to exercise diffing against a sheet that now starts
# with a blank row/column.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# --- Insert a brand-new "Sheet2" right after Sheet1 -----------------------
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet2"

$newSheet.Range("A1").Value = "NewA1"
$newSheet.Range("B1").Value = "NewB1"
$newSheet.Range("C1").Value = "OldC1"

$newSheet.Range("A2").Value = "NewA2"
$newSheet.Range("B2").Value = "NewB2"

$newSheet.Range("A3").Value = "NewA3"
$newSheet.Range("B3").Value = "NewB3"

$newSheet.Range("A4").Value = "NewA4"
$newSheet.Range("B4").Value = "NewB4"

$newSheet.Range("A5").Value = "NewA5"
$newSheet.Range("B5").Value = "NewB5"

$newSheet.Range("A6").Value = "NewA6"
$newSheet.Range("B6").Value = "NewB6"

$newSheet.Range("A7").Value = "NewA7"
$newSheet.Range("B7").Value = "NewB7"

$newSheet.Range("A8").Value = "NewA8"
$newSheet.Range("B8").Value = "NewB8"

$newSheet.Range("A9").Value = "NewA9"
$newSheet.Range("B9").Value = "NewB9"

$newSheet.Range("A10").Value = "NewA10"
$newSheet.Range("B10").Value = "NewB10"

$newSheet.Range("B11").Value = "OldB11"
$newSheet.Range("C11").Value = "OldC11"

# Target stored column width (OOXML "width" attribute) is ~11.855 characters;
# ColumnWidth is offset from the stored width by Excel's fixed ~0.833 padding.
$newSheet.Columns("A").ColumnWidth = 11.02
$newSheet.Range("C12").Select() | Out-Null

# --- Shift Sheet3's data down one row and right one column ---------------
# (insert a blank row 1 and a blank column A so the data now starts at B2).
# NOTE: re-fetch the "Sheet3" reference now that the sheet collection has
# shifted because of the newly inserted "Sheet2" above.
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Rows("1:1").Insert()
$ws3.Columns("A:A").Insert()

# --- Make Sheet3 the active sheet again, with the top row selected -------
$ws3.Activate()
$ws3.Range("A1:XFD1").Select() | Out-Null
